$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header text: bump volume/issue number and report week date range
# ------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/8/2024  Through  7/14/2024"

# ------------------------------------------------------------------
# Weekly crime-stat table refresh (new week of data collected)
# ------------------------------------------------------------------
# Row 15
$ws.Range("F15").Value = 1
$ws.Range("N15").Value = 16.666666666666

# Row 16
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = -13.559322033898
$ws.Range("L16").Value = -30.136986301369
$ws.Range("M16").Value = -32
$ws.Range("N16").Value = -84.307692307692

# Row 17
$ws.Range("I16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 66
$ws.Range("J17").Value = 54
$ws.Range("K17").Value = 22.222222222222
$ws.Range("L17").Value = -8.333333333333
$ws.Range("M17").Value = 60.975609756097
$ws.Range("N17").Value = -49.230769230769

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 79
$ws.Range("J18").Value = 103
$ws.Range("K18").Value = -23.300970873786
$ws.Range("L18").Value = -36.8
$ws.Range("M18").Value = -33.613445378151
$ws.Range("N18").Value = -86.789297658862

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -68.421052631578
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 67
$ws.Range("H19").Value = -28.358208955223
$ws.Range("I19").Value = 334
$ws.Range("J19").Value = 362
$ws.Range("K19").Value = -7.734806629834
$ws.Range("L19").Value = 10.596026490066
$ws.Range("M19").Value = 125.675675675676
$ws.Range("N19").Value = 69.543147208121

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("I16").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 2
$ws.Range("K16").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 52
$ws.Range("J20").Value = 86
$ws.Range("K20").Value = -39.534883720930
$ws.Range("L20").Value = -40.909090909090
$ws.Range("M20").Value = -32.467532467532
$ws.Range("N20").Value = -89.075630252100

# Row 21
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -62.5
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -21.929824561403
$ws.Range("I21").Value = 589
$ws.Range("J21").Value = 667
$ws.Range("K21").Value = -11.694152923538
$ws.Range("L21").Value = -11.295180722891
$ws.Range("M21").Value = 27.765726681128
$ws.Range("N21").Value = -66.051873198847

# Row 22
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

# Row 23
$ws.Range("C14").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("L23").Value = 7.692307692307

# Row 24
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 19.047619047619
$ws.Range("F24").Value = 83
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = 1.219512195121
$ws.Range("I24").Value = 525
$ws.Range("J24").Value = 496
$ws.Range("K24").Value = 5.846774193548
$ws.Range("L24").Value = -3.492647058823
$ws.Range("M24").Value = 82.291666666666

# Row 25
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 58.333333333333
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 21.951219512195
$ws.Range("I25").Value = 300
$ws.Range("J25").Value = 261
$ws.Range("K25").Value = 14.942528735632
$ws.Range("L25").Value = -7.407407407407

# Row 26
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 28.571428571428
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 24.137931034482
$ws.Range("I26").Value = 152
$ws.Range("J26").Value = 128
$ws.Range("K26").Value = 18.75
$ws.Range("L26").Value = 1.333333333333
$ws.Range("M26").Value = 29.914529914529

# Row 27
$ws.Range("F27").Value = 1

# Row 28
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 6
$ws.Range("I28").Value = 17
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = -32
$ws.Range("L28").Value = 21.428571428571

# Row 29
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)

# Row 30
$ws.Range("C14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C30").PasteSpecial(-4122)

# Row 31
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 3
$ws.Range("K31").Value = -50

# Row 33
$ws.Range("C14").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H33").PasteSpecial(-4122)

$ws.Range("A1").Select()
